$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 720, shifting the
# existing rows 720:815 down to 722:817.
$ws.Rows.Item(720).Resize(2).Insert()

# New row 720 (Asterix, 1a (cosecha lavada))
$ws.Cells.Item(720,1).Value = 5
$ws.Cells.Item(720,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(720,3).Value = "Maule"
$ws.Cells.Item(720,4).Value = 45077
$ws.Cells.Item(720,5).Value = 7
$ws.Cells.Item(720,6).Value = 100114001
$ws.Cells.Item(720,7).Value = "Papa"
$ws.Cells.Item(720,8).Value = "Asterix"
$ws.Cells.Item(720,9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(720,10).Value = 1600
$ws.Cells.Item(720,11).Value = 12000
$ws.Cells.Item(720,12).Value = 12000
$ws.Cells.Item(720,13).Value = 12000
$ws.Cells.Item(720,14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(720,15).Value = "Región de Los Lagos"
$ws.Cells.Item(720,16).Value = 480
$ws.Cells.Item(720,17).Value = 25
$ws.Cells.Item(720,18).Value = "Hortaliza"

# New row 721 (Patagonia, 1a (cosecha))
$ws.Cells.Item(721,1).Value = 5
$ws.Cells.Item(721,2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(721,3).Value = "Maule"
$ws.Cells.Item(721,4).Value = 45077
$ws.Cells.Item(721,5).Value = 7
$ws.Cells.Item(721,6).Value = 100114001
$ws.Cells.Item(721,7).Value = "Papa"
$ws.Cells.Item(721,8).Value = "Patagonia"
$ws.Cells.Item(721,9).Value = "1a (cosecha)"
$ws.Cells.Item(721,10).Value = 1600
$ws.Cells.Item(721,11).Value = 10000
$ws.Cells.Item(721,12).Value = 10000
$ws.Cells.Item(721,13).Value = 10000
$ws.Cells.Item(721,14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(721,15).Value = "Región de Los Lagos"
$ws.Cells.Item(721,16).Value = 400
$ws.Cells.Item(721,17).Value = 25
$ws.Cells.Item(721,18).Value = "Hortaliza"
